# Improved modem sleep mode. Fixed #8. Updated BOM.
# Adds three new BOM rows (M5x20 screw, M5x10 screw, M5 hex nut) to the
# "LolinD32-ESP32-SIM800L" sheet, right after the "18650 2000mAh Battery" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LolinD32-ESP32-SIM800L")

$nl = [char]10

# ------------------------------------------------------------------
# 1) Remember every existing hyperlink (its cell, target url and the
#    cached "display" text) so we can re-create them after the new
#    rows shift everything below row 14 three rows further down.
# ------------------------------------------------------------------
$hyperlinkInfo = @()
foreach ($hl in $ws.Hyperlinks) {
    $row = $hl.Range.Row
    $col = $hl.Range.Column
    $hyperlinkInfo += , @($row, $col, $hl.Address, $hl.SubAddress, $hl.TextToDisplay)
}
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2) Insert 3 new blank rows before row 15 (pushes everything from
#    row 15 downward by 3 rows).
# ------------------------------------------------------------------
$ws.Rows("15:17").Insert()

# ---- Row 15 : M5x20 Hexagon socket Head Cap Screw ----
$ws.Rows(15).RowHeight = 45

$ws.Range("A15").Value = 14

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").WrapText = $true
$ws.Range("B15").Value = "M5x20 Hexagon socket Head Cap Screw"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").WrapText = $true

$ws.Range("D15").Value = 2

$ws.Range("E15").Value = "pcs"

$ws.Range("F15").WrapText = $true
$ws.Range("F15").Value = "ISO 4762, A2 Stainless steel, Hexagon/Torx socket head cap screw, k=5mm, l=20mm" + $nl + "http://www.fasteners.eu/standards/ISO/4762"

$ws.Range("G15").WrapText = $true

# ---- Row 16 : M5x10 Hexagon socket Head Cap Screw ----
$ws.Rows(16).RowHeight = 45

$ws.Range("A16").Value = 15

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").Value = "M5x10 Hexagon socket Head Cap Screw"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").WrapText = $true

$ws.Range("D16").Value = 2

$ws.Range("E16").Value = "pcs"

$ws.Range("F16").WrapText = $true
$ws.Range("F16").Value = "ISO 4762, A2  Stainless steel, Hexagon/Torx socket head cap screw, k=5mm, l=10mm" + $nl + "http://www.fasteners.eu/standards/ISO/4762"

$ws.Range("G16").WrapText = $true

# ---- Row 17 : M5 Hexagon Nut ----
$ws.Rows(17).RowHeight = 30

$ws.Range("A17").Value = 16

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").WrapText = $true
$ws.Range("B17").Value = "M5 Hexagon Nut"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").WrapText = $true

$ws.Range("D17").Value = 4

$ws.Range("E17").Value = "pcs"

$ws.Range("F17").WrapText = $true
$ws.Range("F17").Value = "ISO 4032, A2 Stainless steel Hexagon nut, m=~4.5, s=8" + $nl + "http://www.fasteners.eu/standards/ISO/4032"

$ws.Range("G17").WrapText = $true

# ------------------------------------------------------------------
# 3) Renumber the "ID" column (A) for every row that got shifted down.
# ------------------------------------------------------------------
for ($r = 18; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ------------------------------------------------------------------
# 4) Re-create the hyperlinks, shifting anything at/below row 15 down
#    by 3 rows, and restore the original cell text afterwards (Add
#    would otherwise overwrite it with the display text).
# ------------------------------------------------------------------
foreach ($item in $hyperlinkInfo) {
    $origRow = $item[0]
    $col = $item[1]
    $address = $item[2]
    $subAddress = $item[3]
    $display = $item[4]

    $newRow = $origRow
    if ($origRow -ge 15) {
        $newRow = $origRow + 3
    }

    $cell = $ws.Cells.Item($newRow, $col)
    $originalValue = $cell.Value

    $ws.Hyperlinks.Add($cell, $address, $subAddress, "", $display) | Out-Null

    $cell.Value = $originalValue
}

# ------------------------------------------------------------------
# 5) Update the sheet view to match the author's saved selection/scroll
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2:A29").Select()
$excel.ActiveWindow.ScrollRow = 9
